$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "1904"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 781242

$ws.Range("E17").Value = "1904"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 781242

$ws.Range("E18").Value = "1905"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 781242

$ws.Range("E19").Value = "1906"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242

$ws.Range("E20").Value = "1907"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 781242

$ws.Range("E21").Value = "1908"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 781242

$ws.Range("E22").Value = "1909"
$ws.Range("F22").Value = 31249
$ws.Range("G22").Value = 781242

$ws.Range("E23").Value = "1910"
$ws.Range("F23").Value = 21874
$ws.Range("G23").Value = 781242
